# Added road cutback check:
# Insert a new "Checks if road features have any geometry cutbacks" bullet
# immediately before the existing "Geocodes addresses against the road
# centerline data" bullet in the list of validation checks.

$d = $word.ActiveDocument

# Search using a dedicated Range object (not $d.Content directly), since
# Find.Execute collapses *that* range to the match location.
$r = $d.Content
$found = $r.Find.Execute("Checks if all features are inside authoritative boundary", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    $anchorPara = $r.Paragraphs(1)

    # Insert a new paragraph right after the anchor paragraph; it inherits
    # the same paragraph formatting (ListParagraph style / numPr) as the
    # anchor, matching the existing bullet list.
    $anchorPara.Range.InsertParagraphAfter()

    $newPara = $anchorPara.Next()
    $newPara.Range.Text = "Checks if road features have any geometry cutbacks"
}
